# GPLIM-2957: fix excel headers so they match with values in Decision
# https://labopsconfluence.broadinstitute.org/pages/viewpage.action?pageId=22676493
#
# The "Buick Example" sheet's header row used names that didn't line up
# with the field names Decision expects. Rename the two offending
# headers: "Sample ID" -> "Specimen_Number" and "T/N" -> "SAMPLE_TYPE".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Buick Example")
$ws.Activate()

# A1 already carries the "quote prefix" cell style (s="1"), so re-enter
# it with a leading apostrophe to keep that exact style instead of
# Excel minting a brand new (near-duplicate) cell style for plain text.
# F1's style (s="4") has no quote prefix, so a plain value keeps it as-is.
$ws.Range("A1").Formula = "'Specimen_Number"
$ws.Range("F1").Value = "SAMPLE_TYPE"

# Reflect the edited header cells in the current selection.
$ws.Range("F1").Select()
$ws.Range("A1").Select()
